# Update Majuro collection records:
#  - add "latitude" / "longitude" columns (D, E)
#  - fix location-name typos ("Majuru" -> "Majuro", and the second site's
#    name/spelling entirely) for all rows
#  - populate latitude/longitude values for every collection row
#  - widen column C to fit the corrected (longer) location text
#  - leave the selection on B2, matching the authored workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("D1").Value = "latitude"
$ws.Range("E1").Value = "longitude"

# Corrected location names (typo fixes)
$ws.Range("C2").Value = "Airport, Akia, Majuro, Republic of the Marshall Islands"
$ws.Range("C3").Value = "Lojemwe Weto, Ajeltake, Majuro, Republic of the Marshall Islands"
$ws.Range("C4").Value = "Lojemwe Weto, Ajeltake, Majuro, Republic of the Marshall Islands"
$ws.Range("C5").Value = "Lojemwe Weto, Ajeltake, Majuro, Republic of the Marshall Islands"
$ws.Range("C6").Value = "Lojemwe Weto, Ajeltake, Majuro, Republic of the Marshall Islands"
$ws.Range("C7").Value = "Lojemwe Weto, Ajeltake, Majuro, Republic of the Marshall Islands"

# New latitude / longitude data
$ws.Range("D2").Value = 7.068092
$ws.Range("E2").Value = 171.281427

$ws.Range("D3").Value = 7.084167
$ws.Range("E3").Value = 171.133889
$ws.Range("D4").Value = 7.084167
$ws.Range("E4").Value = 171.133889
$ws.Range("D5").Value = 7.084167
$ws.Range("E5").Value = 171.133889
$ws.Range("D6").Value = 7.084167
$ws.Range("E6").Value = 171.133889
$ws.Range("D7").Value = 7.084167
$ws.Range("E7").Value = 171.133889

# Widen the location column to fit the updated (longer) text
$ws.Columns.Item(3).ColumnWidth = 56.5

# Match the saved selection from the authored workbook
$ws.Range("B2").Select() | Out-Null
